$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H7").Value = 11166.5
$ws.Range("J7").Value = 11166.5
$ws.Range("L7").Value = 11166.5
$ws.Range("N7").Value = -11390.5
$ws.Range("H10").Value = 1760
$ws.Range("I10").Value = 600
$ws.Range("J10").Value = 3500
$ws.Range("K10").Value = 600
$ws.Range("L10").Value = 3500
$ws.Range("M10").Value = -307
$ws.Range("N10").Value = -4086
$ws.Range("H14").Value = 11166.5
$ws.Range("J14").Value = 11166.5
$ws.Range("L14").Value = 11166.5
$ws.Range("N14").Value = -11548.5
$ws.Range("H17").Value = 288050.78
$ws.Range("J17").Value = 347282.06
$ws.Range("L17").Value = 1041846.18
$ws.Range("N17").Value = -1042182.18
$ws.Range("H28").Value = 460.22223
$ws.Range("I28").Value = 480.875
$ws.Range("K28").Value = 480.875
$ws.Range("M28").Value = 4.125
$ws.Range("H32").Value = 1180.7142
$ws.Range("I32").Value = 998
$ws.Range("K32").Value = 998
$ws.Range("M32").Value = -672
$ws.Range("H33").Value = 712.2857
$ws.Range("I33").Value = 411.33334
$ws.Range("J33").Value = 938
$ws.Range("K33").Value = 411.33334
$ws.Range("L33").Value = 938
$ws.Range("M33").Value = -182.33334
$ws.Range("N33").Value = -1396
$ws.Range("H38").Value = 5560.7856
$ws.Range("I38").Value = 232
$ws.Range("J38").Value = 12665.833
$ws.Range("K38").Value = 696
$ws.Range("L38").Value = 37997.499
$ws.Range("M38").Value = -324
$ws.Range("N38").Value = -38741.499
$ws.Range("H39").Value = 1037.909
$ws.Range("I39").Value = 888.1667
$ws.Range("K39").Value = 2664.5001
$ws.Range("M39").Value = -2368.5001
$ws.Range("H40").Value = 10347512
$ws.Range("I40").Value = 3076.8262
$ws.Range("J40").Value = 50001180
$ws.Range("K40").Value = 3076.8262
$ws.Range("L40").Value = 50001180
$ws.Range("M40").Value = -2901.8262
$ws.Range("N40").Value = -50001530
$ws.Range("H41").Value = 1600.7858
$ws.Range("I41").Value = 616.8333
$ws.Range("J41").Value = 2338.75
$ws.Range("K41").Value = 616.8333
$ws.Range("L41").Value = 2338.75
$ws.Range("M41").Value = -176.8333
$ws.Range("N41").Value = -3218.75
$ws.Range("H42").Value = 259.3846
$ws.Range("J42").Value = 316.66666
$ws.Range("L42").Value = 949.9999799999999
$ws.Range("N42").Value = -1409.99998
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H47").Value = 29800
$ws.Range("I47").Value = 29800
$ws.Range("K47").Value = 29800
$ws.Range("M47").Value = -28828
$ws.Range("H48").Value = 1999
$ws.Range("J48").Value = 1999
$ws.Range("L48").Value = 5997
$ws.Range("N48").Value = -6581
$ws.Range("H51").Value = 22200
$ws.Range("I51").Value = 23667
$ws.Range("J51").Value = 19999.5
$ws.Range("K51").Value = 23667
$ws.Range("L51").Value = 19999.5
$ws.Range("M51").Value = -23183
$ws.Range("N51").Value = -20967.5
$ws.Range("H56").Value = 1999
$ws.Range("J56").Value = 1999
$ws.Range("L56").Value = 5997
$ws.Range("N56").Value = -7065
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H62").Value = 68187840
$ws.Range("I62").Value = 35720372
$ws.Range("J62").Value = 125005910
$ws.Range("K62").Value = 35720372
$ws.Range("L62").Value = 125005910
$ws.Range("M62").Value = -35719748
$ws.Range("N62").Value = -125007158
$ws.Range("H65").Value = 68187840
$ws.Range("I65").Value = 35720372
$ws.Range("J65").Value = 125005910
$ws.Range("K65").Value = 178601860
$ws.Range("L65").Value = 625029550
$ws.Range("M65").Value = -178598740
$ws.Range("N65").Value = -625035790
$ws.Range("H98").Value = 3739052.2
$ws.Range("I98").Value = 3953760
$ws.Range("K98").Value = 3953760
$ws.Range("M98").Value = -3952262
$ws.Range("H115").Value = 81779896
$ws.Range("J115").Value = 9999.5
$ws.Range("L115").Value = 29998.5
$ws.Range("N115").Value = -33132.5
$ws.Range("H122").Value = 3739052.2
$ws.Range("I122").Value = 3953760
$ws.Range("K122").Value = 11861280
$ws.Range("M122").Value = -11858830
$ws.Range("H132").Value = 1460.7407
$ws.Range("I132").Value = 1137.0488
$ws.Range("J132").Value = 2481.6155
$ws.Range("K132").Value = 3411.1464
$ws.Range("L132").Value = 7444.8465
$ws.Range("M132").Value = -881.1464000000001
$ws.Range("N132").Value = -12504.8465
$ws.Range("H138").Value = 5345.794
$ws.Range("J138").Value = 5873.6787
$ws.Range("L138").Value = 17621.0361
$ws.Range("N138").Value = -27901.0361

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7262.9155
$ws.Range("I32").Value = 7262.9155
$ws.Range("K32").Value = 7262.9155
$ws.Range("M32").Value = -6975.9155
$ws.Range("H43").Value = 24996.666
$ws.Range("J43").Value = 24996.666
$ws.Range("L43").Value = 24996.666
$ws.Range("N43").Value = -25622.666
$ws.Range("H61").Value = 5216.9414
$ws.Range("I61").Value = 2320.889
$ws.Range("K61").Value = 2320.889
$ws.Range("M61").Value = -2108.889
$ws.Range("H74").Value = 134572.8
$ws.Range("J74").Value = 4014
$ws.Range("L74").Value = 4014
$ws.Range("N74").Value = -5762
$ws.Range("H77").Value = 134572.8
$ws.Range("J77").Value = 4014
$ws.Range("L77").Value = 20070
$ws.Range("N77").Value = -28806
$ws.Range("H97").Value = 908.375
$ws.Range("I97").Value = 595.4
$ws.Range("K97").Value = 595.4
$ws.Range("M97").Value = -99.39999999999998
$ws.Range("H102").Value = 1942.7354
$ws.Range("I102").Value = 1773.5
$ws.Range("J102").Value = 2492.75
$ws.Range("K102").Value = 1773.5
$ws.Range("L102").Value = 2492.75
$ws.Range("M102").Value = -151.5
$ws.Range("N102").Value = -5736.75
$ws.Range("H122").Value = 3098.7778
$ws.Range("I122").Value = 3055.3333
$ws.Range("J122").Value = 3316
$ws.Range("K122").Value = 9165.999899999999
$ws.Range("L122").Value = 9948
$ws.Range("M122").Value = -6715.999899999999
$ws.Range("N122").Value = -14848
$ws.Range("H132").Value = 3162.2812
$ws.Range("I132").Value = 3818.9583
$ws.Range("J132").Value = 1192.25
$ws.Range("K132").Value = 11456.8749
$ws.Range("L132").Value = 3576.75
$ws.Range("M132").Value = -8926.874899999999
$ws.Range("N132").Value = -8636.75
$ws.Range("H136").Value = 5216.9414
$ws.Range("I136").Value = 2320.889
$ws.Range("K136").Value = 6962.667
$ws.Range("M136").Value = -4412.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 44942.5
$ws.Range("J46").Value = 69885
$ws.Range("L46").Value = 69885
$ws.Range("N46").Value = -70481
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H54").Value = 1665
$ws.Range("I54").Value = 1665
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 1665
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -1181
$ws.Range("N54").ClearContents()
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H86").Value = 2132.7144
$ws.Range("J86").Value = 2999
$ws.Range("L86").Value = 2999
$ws.Range("N86").Value = -5245
$ws.Range("H89").Value = 2132.7144
$ws.Range("J89").Value = 2999
$ws.Range("L89").Value = 14995
$ws.Range("N89").Value = -26227
$ws.Range("H132").Value = 114998
$ws.Range("J132").Value = 114998
$ws.Range("L132").Value = 114998
$ws.Range("N132").Value = -125118

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 200005800
$ws.Range("J4").Value = 200005800
$ws.Range("L4").Value = 200005800
$ws.Range("N4").Value = -200006024
$ws.Range("H31").Value = 201621.56
$ws.Range("J31").Value = 50638.875
$ws.Range("L31").Value = 50638.875
$ws.Range("N31").Value = -51228.875
$ws.Range("H32").Value = 7369.7
$ws.Range("I32").Value = 5966.3335
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 5966.3335
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -5650.3335
$ws.Range("N32").Value = -20632
$ws.Range("H34").Value = 201621.56
$ws.Range("J34").Value = 50638.875
$ws.Range("L34").Value = 50638.875
$ws.Range("N34").Value = -51042.875
$ws.Range("H55").Value = 7326.6665
$ws.Range("I55").Value = 6499.5
$ws.Range("J55").Value = 8981
$ws.Range("K55").Value = 6499.5
$ws.Range("L55").Value = 8981
$ws.Range("M55").Value = -6184.5
$ws.Range("N55").Value = -9611
$ws.Range("H58").Value = 1956.4762
$ws.Range("I58").Value = 1823.4667
$ws.Range("K58").Value = 1823.4667
$ws.Range("M58").Value = -1620.4667
$ws.Range("H69").Value = 20499.5
$ws.Range("I69").Value = 20499.5
$ws.Range("K69").Value = 20499.5
$ws.Range("M69").Value = -19750.5
$ws.Range("H72").Value = 20499.5
$ws.Range("I72").Value = 20499.5
$ws.Range("K72").Value = 61498.5
$ws.Range("M72").Value = -57754.5
$ws.Range("H86").Value = 3745.0715
$ws.Range("I86").Value = 3638.9
$ws.Range("K86").Value = 3638.9
$ws.Range("M86").Value = -2515.9
$ws.Range("H89").Value = 3745.0715
$ws.Range("I89").Value = 3638.9
$ws.Range("K89").Value = 18194.5
$ws.Range("M89").Value = -12578.5
$ws.Range("H107").Value = 4515.5405
$ws.Range("J107").Value = 6489.5
$ws.Range("L107").Value = 6489.5
$ws.Range("N107").Value = -10329.5
$ws.Range("H122").Value = 752.375
$ws.Range("I122").Value = 752.375
$ws.Range("K122").Value = 2257.125
$ws.Range("M122").Value = 192.875
$ws.Range("H132").Value = 2699.5625
$ws.Range("I132").Value = 2746.8667
$ws.Range("J132").Value = 1990
$ws.Range("K132").Value = 8240.6001
$ws.Range("L132").Value = 5970
$ws.Range("M132").Value = -5710.6001
$ws.Range("N132").Value = -11030
$ws.Range("H134").Value = 5646.25
$ws.Range("I134").Value = 6623.857
$ws.Range("J134").Value = 2713.4285
$ws.Range("K134").Value = 19871.571
$ws.Range("L134").Value = 8140.2855
$ws.Range("M134").Value = -17336.571
$ws.Range("N134").Value = -13210.2855
$ws.Range("H136").Value = 1956.4762
$ws.Range("I136").Value = 1823.4667
$ws.Range("K136").Value = 5470.4001
$ws.Range("M136").Value = -2920.4001
$ws.Range("H139").Value = 69780
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 69780
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 69780
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -80060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 54.5
$ws.Range("I2").Value = 42.307693
$ws.Range("K2").Value = 253.846158
$ws.Range("M2").Value = -140.846158
$ws.Range("H5").Value = 807.8333
$ws.Range("I5").Value = 807.8333
$ws.Range("K5").Value = 2423.4999
$ws.Range("M5").Value = -2311.4999
$ws.Range("H16").Value = 3444.5
$ws.Range("J16").Value = 3889
$ws.Range("L16").Value = 11667
$ws.Range("N16").Value = -12013
$ws.Range("H39").Value = 1625.5227
$ws.Range("I39").Value = 915.6667
$ws.Range("J39").Value = 1737.6052
$ws.Range("K39").Value = 2747.0001
$ws.Range("L39").Value = 5212.8156
$ws.Range("M39").Value = -2453.0001
$ws.Range("N39").Value = -5800.8156
$ws.Range("H55").Value = 7212.5557
$ws.Range("I55").Value = 4223.5
$ws.Range("K55").Value = 12670.5
$ws.Range("M55").Value = -12493.5
$ws.Range("H70").Value = 1714.1428
$ws.Range("I70").Value = 999.5
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 2998.5
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -2683.5
$ws.Range("N70").Value = -6630
$ws.Range("H73").Value = 1714.1428
$ws.Range("I73").Value = 999.5
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 2998.5
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -1906.5
$ws.Range("N73").Value = -8184
$ws.Range("H75").Value = 1786.4
$ws.Range("I75").Value = 1487.2858
$ws.Range("K75").Value = 4461.857400000001
$ws.Range("M75").Value = -3463.857400000001
$ws.Range("H78").Value = 1786.4
$ws.Range("I78").Value = 1487.2858
$ws.Range("K78").Value = 13385.5722
$ws.Range("M78").Value = -8393.572200000001
$ws.Range("H107").Value = 511.04544
$ws.Range("J107").Value = 520
$ws.Range("L107").Value = 1560
$ws.Range("N107").Value = -5400
$ws.Range("H113").Value = 1182.6666
$ws.Range("J113").Value = 499
$ws.Range("L113").Value = 1497
$ws.Range("N113").Value = -5837
$ws.Range("H131").Value = 5683249
$ws.Range("I131").Value = 125000776
$ws.Range("J131").Value = 1462.1548
$ws.Range("K131").Value = 375002328
$ws.Range("L131").Value = 4386.4644
$ws.Range("M131").Value = -374997288
$ws.Range("N131").Value = -14466.4644
$ws.Range("H132").Value = 5415.4644
$ws.Range("I132").Value = 6451.9546
$ws.Range("J132").Value = 1615
$ws.Range("K132").Value = 58067.5914
$ws.Range("L132").Value = 14535
$ws.Range("M132").Value = -55537.5914
$ws.Range("N132").Value = -19595
$ws.Range("H135").Value = 807.8333
$ws.Range("I135").Value = 807.8333
$ws.Range("K135").Value = 7270.4997
$ws.Range("M135").Value = -4735.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7772.625
$ws.Range("I70").Value = 8321.117
$ws.Range("J70").Value = 6440.5713
$ws.Range("K70").Value = 8321.117
$ws.Range("L70").Value = 6440.5713
$ws.Range("M70").Value = -8051.117
$ws.Range("N70").Value = -6980.5713
$ws.Range("H73").Value = 7772.625
$ws.Range("I73").Value = 8321.117
$ws.Range("J73").Value = 6440.5713
$ws.Range("K73").Value = 8321.117
$ws.Range("L73").Value = 6440.5713
$ws.Range("M73").Value = -7385.117
$ws.Range("N73").Value = -8312.5713
$ws.Range("H80").Value = 3680.2
$ws.Range("I80").Value = 3924.48
$ws.Range("K80").Value = 3924.48
$ws.Range("M80").Value = -2926.48
$ws.Range("H83").Value = 3680.2
$ws.Range("I83").Value = 3924.48
$ws.Range("K83").Value = 19622.4
$ws.Range("M83").Value = -14630.4
$ws.Range("H102").Value = 2369.9333
$ws.Range("I102").Value = 2465.4614
$ws.Range("J102").Value = 1749
$ws.Range("K102").Value = 2465.4614
$ws.Range("L102").Value = 1749
$ws.Range("M102").Value = -843.4614000000001
$ws.Range("N102").Value = -4993
$ws.Range("H122").Value = 2105.6875
$ws.Range("I122").Value = 2216.5386
$ws.Range("J122").Value = 1625.3334
$ws.Range("K122").Value = 6649.6158
$ws.Range("L122").Value = 4876.0002
$ws.Range("M122").Value = -4199.6158
$ws.Range("N122").Value = -9776.0002
$ws.Range("H126").Value = 8580.875
$ws.Range("J126").Value = 4219.75
$ws.Range("L126").Value = 12659.25
$ws.Range("N126").Value = -17599.25
$ws.Range("H132").Value = 50536.5
$ws.Range("I132").Value = 61748.668
$ws.Range("J132").Value = 16900
$ws.Range("K132").Value = 185246.004
$ws.Range("L132").Value = 50700
$ws.Range("M132").Value = -182716.004
$ws.Range("N132").Value = -55760

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8417.138999999999
$ws.Range("I7").Value = 8516.25
$ws.Range("K7").Value = 8516.25
$ws.Range("M7").Value = -8404.25
$ws.Range("H16").Value = 622.85
$ws.Range("I16").Value = 682.05554
$ws.Range("J16").Value = 90
$ws.Range("K16").Value = 682.05554
$ws.Range("L16").Value = 90
$ws.Range("M16").Value = -512.05554
$ws.Range("N16").Value = -430
$ws.Range("H22").Value = 1392
$ws.Range("I22").Value = 1340.75
$ws.Range("J22").Value = 1404.8125
$ws.Range("K22").Value = 1340.75
$ws.Range("L22").Value = 1404.8125
$ws.Range("M22").Value = -1045.75
$ws.Range("N22").Value = -1994.8125
$ws.Range("H25").Value = 1000
$ws.Range("I25").Value = 1000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 1000
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -770
$ws.Range("N25").ClearContents()
$ws.Range("H27").Value = 1392
$ws.Range("I27").Value = 1340.75
$ws.Range("J27").Value = 1404.8125
$ws.Range("K27").Value = 1340.75
$ws.Range("L27").Value = 1404.8125
$ws.Range("M27").Value = -1233.75
$ws.Range("N27").Value = -1618.8125
$ws.Range("H40").Value = 8793.6
$ws.Range("I40").Value = 7996.6665
$ws.Range("K40").Value = 7996.6665
$ws.Range("M40").Value = -7860.6665
$ws.Range("H107").Value = 261008.25
$ws.Range("I107").Value = 261008.25
$ws.Range("K107").Value = 261008.25
$ws.Range("M107").Value = -259088.25
$ws.Range("H122").Value = 9155
$ws.Range("I122").Value = 8698.799999999999
$ws.Range("J122").Value = 9915.333000000001
$ws.Range("K122").Value = 26096.4
$ws.Range("L122").Value = 29745.999
$ws.Range("M122").Value = -23646.4
$ws.Range("N122").Value = -34645.999
$ws.Range("H126").Value = 8417.138999999999
$ws.Range("I126").Value = 8516.25
$ws.Range("K126").Value = 25548.75
$ws.Range("M126").Value = -23078.75
$ws.Range("H132").Value = 9995
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 9995
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 29985
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -35045
$ws.Range("H136").Value = 4462.1665
$ws.Range("I136").Value = 4442
$ws.Range("J136").Value = 4472.25
$ws.Range("K136").Value = 13326
$ws.Range("L136").Value = 13416.75
$ws.Range("M136").Value = -10776
$ws.Range("N136").Value = -18516.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 1210.5
$ws.Range("J8").Value = 1210.5
$ws.Range("L8").Value = 1210.5
$ws.Range("N8").Value = -1490.5
$ws.Range("H81").Value = 2115
$ws.Range("I81").Value = 2149.6
$ws.Range("J81").Value = 1596
$ws.Range("K81").Value = 4299.2
$ws.Range("L81").Value = 3192
$ws.Range("M81").Value = -3238.2
$ws.Range("N81").Value = -5314
$ws.Range("H84").Value = 2115
$ws.Range("I84").Value = 2149.6
$ws.Range("J84").Value = 1596
$ws.Range("K84").Value = 21496
$ws.Range("L84").Value = 15960
$ws.Range("M84").Value = -16192
$ws.Range("N84").Value = -26568
$ws.Range("H96").Value = 103097.4
$ws.Range("I96").Value = 127621.875
$ws.Range("K96").Value = 127621.875
$ws.Range("M96").Value = -126248.875
$ws.Range("H113").Value = 1253.762
$ws.Range("I113").Value = 1352.75
$ws.Range("J113").Value = 1121.7778
$ws.Range("K113").Value = 4058.25
$ws.Range("L113").Value = 3365.3334
$ws.Range("M113").Value = -1888.25
$ws.Range("N113").Value = -7705.3334
$ws.Range("H126").Value = 2598.8
$ws.Range("I126").Value = 2498.75
$ws.Range("K126").Value = 7496.25
$ws.Range("M126").Value = -5026.25
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 207363.58
$ws.Range("I136").Value = 281442.8
$ws.Range("K136").Value = 844328.3999999999
$ws.Range("M136").Value = -841778.3999999999
